$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 1 values (B1:E1)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Move the CON value from D2 to C2
$ws.Range("D2").ClearContents()
$ws.Range("C2").Value = 27.73009143525185

# Clear STR values in B3 and C3
$ws.Range("B3").ClearContents()
$ws.Range("C3").ClearContents()

# Update the selected range to match the saved selection state
$ws.Range("B1:E3").Select()
